$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Input")
$ws2 = $wb.Worksheets.Item("Output")

# Sheet1 ("Input") content updates
$ws1.Range("B7").Value = "Jhon Deer"
$ws1.Range("A8").Value = "GroupAddClient"
$ws1.Range("B8").Value = "click"
$ws1.Range("B2").Value = "Group4109"

# Sheet2 ("Output") content updates
$ws2.Range("A2").Value = "VerifyClientCreated"
$ws2.Range("B1").Value = "Group4109"
$ws2.Range("B1").Font.Name = "Calibri"
$ws2.Range("B1").WrapText = $true
$ws2.Range("B2").Value = "Jhon Deer"

# Final selections / active sheet, matching the saved view state
$ws2.Activate() | Out-Null
$ws1.Range("B2").Select() | Out-Null
$ws2.Range("A5").Select() | Out-Null
